$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values (price/volume refresh + 3 new coin-list rows shifted in).
# Using a leading apostrophe forces Excel to store these as literal text (matching the
# source inlineStr cells) instead of inferring numbers/percentages from numeric-looking strings.

$ws.Range("D2").Value = "'" + "243.60"
$ws.Range("E2").Value = "'" + "-0.25%"
$ws.Range("E3").Value = "'" + "13.67%"
$ws.Range("D4").Value = "'" + "5.141"
$ws.Range("E4").Value = "'" + "0.21%"
$ws.Range("D5").Value = "'" + "0.05670"
$ws.Range("E5").Value = "'" + "1.44%"
$ws.Range("D6").Value = "'" + "6.527"
$ws.Range("E6").Value = "'" + "0.80%"
$ws.Range("D7").Value = "'" + "0.8407"
$ws.Range("E7").Value = "'" + "2.45%"
$ws.Range("D8").Value = "'" + "0.8651"
$ws.Range("E8").Value = "'" + "3.49%"
$ws.Range("B9").Value = "'" + "One"
$ws.Range("C9").Value = "'" + "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'" + "0.0005986"
$ws.Range("E9").Value = "'" + "-0.07%"
$ws.Range("B10").Value = "'" + "WazirX"
$ws.Range("C10").Value = "'" + "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'" + "0.1337"
$ws.Range("E10").Value = "'" + "0.41%"
$ws.Range("B11").Value = "'" + "MandalaExchangeToken"
$ws.Range("C11").Value = "'" + "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'" + "0.06915"
$ws.Range("E11").Value = "'" + "-1.03%"
$ws.Range("B12").Value = "'" + "BitrueCoin"
$ws.Range("C12").Value = "'" + "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'" + "0.02903"
$ws.Range("E12").Value = "'" + "0.53%"
$ws.Range("B13").Value = "'" + "BitMartToken"
$ws.Range("C13").Value = "'" + "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'" + "0.09377"
$ws.Range("E13").Value = "'" + "-0.07%"
$ws.Range("B14").Value = "'" + "BitForexToken"
$ws.Range("C14").Value = "'" + "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'" + "0.001513"
$ws.Range("E14").Value = "'" + "0.04%"
$ws.Range("B15").Value = "'" + "CoinExToken"
$ws.Range("C15").Value = "'" + "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'" + "0.04161"
$ws.Range("E15").Value = "'" + "-10.23%"
$ws.Range("D16").Value = "'" + "0.006195"
$ws.Range("E16").Value = "'" + "-0.89%"
$ws.Range("D17").Value = "'" + "3.507"
$ws.Range("E17").Value = "'" + "-3.93%"
$ws.Range("E18").Value = "'" + "-0.31%"
$ws.Range("D21").Value = "'" + "0.03261"
$ws.Range("E21").Value = "'" + "4.55%"
$ws.Range("D22").Value = "'" + "0.1296"
$ws.Range("E22").Value = "'" + "-0.27%"
$ws.Range("D23").Value = "'" + "3.616"
$ws.Range("E23").Value = "'" + "-3.40%"
$ws.Range("D24").Value = "'" + "0.1373"
$ws.Range("E24").Value = "'" + "-0.07%"
$ws.Range("D25").Value = "'" + "0.001210"
$ws.Range("E25").Value = "'" + "-2.87%"
$ws.Range("D26").Value = "'" + "0.004442"
$ws.Range("E26").Value = "'" + "-1.24%"
$ws.Range("D27").Value = "'" + "0.0001179"
$ws.Range("E27").Value = "'" + "22.87%"
$ws.Range("E28").Value = "'" + "0.26%"
$ws.Range("D40").Value = "'" + "0.03711"
$ws.Range("E40").Value = "'" + "1.99%"
$ws.Range("D41").Value = "'" + "0.005327"
$ws.Range("E41").Value = "'" + "55.40%"
$ws.Range("D42").Value = "'" + "0.1056"
$ws.Range("E42").Value = "'" + "-22.52%"
$ws.Range("D43").Value = "'" + "0.002310"
$ws.Range("E43").Value = "'" + "-11.83%"
$ws.Range("D44").Value = "'" + "0.009812"
$ws.Range("E44").Value = "'" + "10.79%"
$ws.Range("D45").Value = "'" + "0.00005108"
$ws.Range("E45").Value = "'" + "-4.42%"
$ws.Range("E46").Value = "'" + "-0.08%"
$ws.Range("D47").Value = "'" + "0.09994"
$ws.Range("E47").Value = "'" + "-30.60%"
$ws.Range("E48").Value = "'" + "18.70%"
$ws.Range("D49").Value = "'" + "0.00002099"
$ws.Range("E49").Value = "'" + "-0.08%"
$ws.Range("D50").Value = "'" + "0.0001999"
$ws.Range("E50").Value = "'" + "-0.08%"
